$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 5000
$ws.Range("I13").Value = 5000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -4831
$ws.Range("N13").ClearContents()
$ws.Range("H16").Value = 19997.5
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 19997.5
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 19997.5
$ws.Range("N16").Value = -20457.5
$ws.Range("H17").Value = 144542.72
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 144542.72
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 433628.16
$ws.Range("N17").Value = -433964.16
$ws.Range("H33").Value = 198.41176
$ws.Range("I33").Value = 205.08333
$ws.Range("J33").Value = 182.4
$ws.Range("K33").Value = 205.08333
$ws.Range("L33").Value = 182.4
$ws.Range("M33").Value = 23.91667000000001
$ws.Range("N33").Value = -640.4
$ws.Range("H38").Value = 2435.5
$ws.Range("I38").Value = 1706.1111
$ws.Range("J38").Value = 9000
$ws.Range("K38").Value = 5118.3333
$ws.Range("L38").Value = 27000
$ws.Range("M38").Value = -4746.3333
$ws.Range("N38").Value = -27744
$ws.Range("H40").Value = 3630.5264
$ws.Range("I40").Value = 2855.5
$ws.Range("J40").Value = 5800.6
$ws.Range("K40").Value = 2855.5
$ws.Range("L40").Value = 5800.6
$ws.Range("M40").Value = -2680.5
$ws.Range("N40").Value = -6150.6
$ws.Range("H42").Value = 4306.8
$ws.Range("I42").Value = 178.66667
$ws.Range("J42").Value = 10499
$ws.Range("K42").Value = 536.00001
$ws.Range("L42").Value = 31497
$ws.Range("M42").Value = -306.00001
$ws.Range("N42").Value = -31957
$ws.Range("H74").Value = 23861766
$ws.Range("I74").Value = 23861766
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 23861766
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -23860830
$ws.Range("H76").Value = 7981
$ws.Range("I76").Value = 7439.4
$ws.Range("J76").Value = 8883.666999999999
$ws.Range("K76").Value = 7439.4
$ws.Range("L76").Value = 8883.666999999999
$ws.Range("M76").Value = -7124.4
$ws.Range("N76").Value = -9513.666999999999
$ws.Range("H77").Value = 23861766
$ws.Range("I77").Value = 23861766
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 119308830
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -119304150
$ws.Range("H79").Value = 7981
$ws.Range("I79").Value = 7439.4
$ws.Range("J79").Value = 8883.666999999999
$ws.Range("K79").Value = 7439.4
$ws.Range("L79").Value = 8883.666999999999
$ws.Range("M79").Value = -6347.4
$ws.Range("N79").Value = -11067.667
$ws.Range("H98").Value = 1577.6923
$ws.Range("I98").Value = 1626.6666
$ws.Range("J98").Value = 990
$ws.Range("K98").Value = 1626.6666
$ws.Range("L98").Value = 990
$ws.Range("M98").Value = -128.6666
$ws.Range("N98").Value = -3986
$ws.Range("H122").Value = 1577.6923
$ws.Range("I122").Value = 1626.6666
$ws.Range("J122").Value = 990
$ws.Range("K122").Value = 4879.9998
$ws.Range("L122").Value = 2970
$ws.Range("M122").Value = -2429.9998
$ws.Range("N122").Value = -7870
$ws.Range("H125").Value = 1660.4286
$ws.Range("I125").Value = 1597.8334
$ws.Range("J125").Value = 2036
$ws.Range("K125").Value = 14380.5006
$ws.Range("L125").Value = 18324
$ws.Range("M125").Value = -11920.5006
$ws.Range("N125").Value = -23244
$ws.Range("H129").Value = 1969.5
$ws.Range("I129").Value = 1941.3334
$ws.Range("J129").Value = 1997.6666
$ws.Range("K129").Value = 5824.0002
$ws.Range("L129").Value = 5992.9998
$ws.Range("M129").Value = -824.0002000000004
$ws.Range("N129").Value = -15992.9998
$ws.Range("H132").Value = 2318
$ws.Range("I132").Value = 2650.8125
$ws.Range("J132").Value = 1430.5
$ws.Range("K132").Value = 7952.4375
$ws.Range("L132").Value = 4291.5
$ws.Range("M132").Value = -5422.4375
$ws.Range("N132").Value = -9351.5
$ws.Range("H138").Value = 2519.6592
$ws.Range("I138").Value = 2298.2778
$ws.Range("J138").Value = 2672.923
$ws.Range("K138").Value = 6894.8334
$ws.Range("L138").Value = 8018.768999999999
$ws.Range("M138").Value = -1754.8334
$ws.Range("N138").Value = -18298.769

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H30").Value = 300
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 300
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 300
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -600
$ws.Range("H63").Value = 3980
$ws.Range("I63").Value = 4000
$ws.Range("J63").Value = 3950
$ws.Range("K63").Value = 4000
$ws.Range("L63").Value = 3950
$ws.Range("M63").Value = -3314
$ws.Range("N63").Value = -5322
$ws.Range("H66").Value = 3980
$ws.Range("I66").Value = 4000
$ws.Range("J66").Value = 3950
$ws.Range("K66").Value = 20000
$ws.Range("L66").Value = 19750
$ws.Range("M66").Value = -16568
$ws.Range("N66").Value = -26614
$ws.Range("H114").Value = 104799
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 104799
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 104799
$ws.Range("N114").Value = -113477
$ws.Range("H117").Value = 34999
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 34999
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 34999
$ws.Range("N117").Value = -44177
$ws.Range("H121").Value = 24231
$ws.Range("I121").Value = 24231
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 24231
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -22484
$ws.Range("N121").ClearContents()
$ws.Range("H123").Value = 69248
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 69248
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 69248
$ws.Range("N123").Value = -79048

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 20300
$ws.Range("I51").Value = 20300
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 20300
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -19809
$ws.Range("H53").Value = 87780
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 87780
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 87780
$ws.Range("N53").Value = -88928

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 1262.5
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1262.5
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1262.5
$ws.Range("N13").Value = -1540.5
$ws.Range("H86").Value = 3852.889
$ws.Range("I86").Value = 4556.3335
$ws.Range("J86").Value = 3501.1667
$ws.Range("K86").Value = 4556.3335
$ws.Range("L86").Value = 3501.1667
$ws.Range("M86").Value = -3433.3335
$ws.Range("N86").Value = -5747.1667
$ws.Range("H89").Value = 3852.889
$ws.Range("I89").Value = 4556.3335
$ws.Range("J89").Value = 3501.1667
$ws.Range("K89").Value = 22781.6675
$ws.Range("L89").Value = 17505.8335
$ws.Range("M89").Value = -17165.6675
$ws.Range("N89").Value = -28737.8335
$ws.Range("H122").Value = 1828.9
$ws.Range("I122").Value = 1580.9642
$ws.Range("J122").Value = 5300
$ws.Range("K122").Value = 4742.892599999999
$ws.Range("L122").Value = 15900
$ws.Range("M122").Value = -2292.892599999999
$ws.Range("N122").Value = -20800
$ws.Range("H132").Value = 23811020
$ws.Range("I132").Value = 25001460
$ws.Range("J132").Value = 2204.5
$ws.Range("K132").Value = 75004380
$ws.Range("L132").Value = 6613.5
$ws.Range("M132").Value = -75001850
$ws.Range("N132").Value = -11673.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2586500
$ws.Range("I4").Value = 1574214.1
$ws.Range("J4").Value = 4003700
$ws.Range("K4").Value = 4722642.300000001
$ws.Range("L4").Value = 12011100
$ws.Range("M4").Value = -4722530.300000001
$ws.Range("N4").Value = -12011324
$ws.Range("H14").Value = 377.44446
$ws.Range("I14").Value = 377.44446
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1132.33338
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -959.33338
$ws.Range("H34").Value = 1898.875
$ws.Range("I34").Value = 459.8889
$ws.Range("J34").Value = 3749
$ws.Range("K34").Value = 1379.6667
$ws.Range("L34").Value = 11247
$ws.Range("M34").Value = -1295.6667
$ws.Range("N34").Value = -11415
$ws.Range("H55").Value = 4500.364
$ws.Range("I55").Value = 2252
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 6756
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = -6579
$ws.Range("N55").Value = -15354

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 2999
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 2999
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 2999
$ws.Range("N54").Value = -3779
$ws.Range("H70").Value = 16899.5
$ws.Range("I70").Value = 16899.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 16899.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -16629.5
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 16899.5
$ws.Range("I73").Value = 16899.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 16899.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -15963.5
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 2232.9443
$ws.Range("I80").Value = 2153.923
$ws.Range("J80").Value = 2438.4
$ws.Range("K80").Value = 2153.923
$ws.Range("L80").Value = 2438.4
$ws.Range("M80").Value = -1155.923
$ws.Range("N80").Value = -4434.4
$ws.Range("H83").Value = 2232.9443
$ws.Range("I83").Value = 2153.923
$ws.Range("J83").Value = 2438.4
$ws.Range("K83").Value = 10769.615
$ws.Range("L83").Value = 12192
$ws.Range("M83").Value = -5777.614999999998
$ws.Range("N83").Value = -22176
$ws.Range("H97").Value = 957.8570999999999
$ws.Range("I97").Value = 605.5
$ws.Range("J97").Value = 1662.5714
$ws.Range("K97").Value = 605.5
$ws.Range("L97").Value = 1662.5714
$ws.Range("M97").Value = -109.5
$ws.Range("N97").Value = -2654.5714
$ws.Range("H103").Value = 19499.5
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 19499.5
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 19499.5
$ws.Range("N103").Value = -21843.5
$ws.Range("H105").Value = 59998
$ws.Range("I105").Value = 54998
$ws.Range("J105").Value = 62498
$ws.Range("K105").Value = 54998
$ws.Range("L105").Value = 62498
$ws.Range("M105").Value = -51504
$ws.Range("N105").Value = -69486
$ws.Range("H122").Value = 7572
$ws.Range("I122").Value = 5252
$ws.Range("J122").Value = 8500
$ws.Range("K122").Value = 15756
$ws.Range("L122").Value = 25500
$ws.Range("M122").Value = -13306
$ws.Range("N122").Value = -30400

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1099.5
$ws.Range("I46").Value = 970.8570999999999
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 970.8570999999999
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -782.8570999999999
$ws.Range("N46").Value = -2376
$ws.Range("H68").Value = 5000978
$ws.Range("I68").Value = 6250950
$ws.Range("J68").Value = 1090
$ws.Range("K68").Value = 6250950
$ws.Range("L68").Value = 1090
$ws.Range("M68").Value = -6250201
$ws.Range("N68").Value = -2588
$ws.Range("H71").Value = 5000978
$ws.Range("I71").Value = 6250950
$ws.Range("J71").Value = 1090
$ws.Range("K71").Value = 31254750
$ws.Range("L71").Value = 5450
$ws.Range("M71").Value = -31251006
$ws.Range("N71").Value = -12938
$ws.Range("H82").Value = 955.625
$ws.Range("I82").Value = 944.3333
$ws.Range("J82").Value = 989.5
$ws.Range("K82").Value = 944.3333
$ws.Range("L82").Value = 989.5
$ws.Range("M82").Value = -583.3333
$ws.Range("N82").Value = -1711.5
$ws.Range("H85").Value = 955.625
$ws.Range("I85").Value = 944.3333
$ws.Range("J85").Value = 989.5
$ws.Range("K85").Value = 944.3333
$ws.Range("L85").Value = 989.5
$ws.Range("M85").Value = 303.6667
$ws.Range("N85").Value = -3485.5
